$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5200
$ws.Range("I40").Value = 6600
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 6600
$ws.Range("L40").Value = 1000
$ws.Range("M40").Value = -6425
$ws.Range("N40").Value = -1350

$ws.Range("H80").Value = 999.875
$ws.Range("I80").Value = 500
$ws.Range("J80").Value = 1071.2858
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 3213.8574
$ws.Range("M80").Value = -502
$ws.Range("N80").Value = -5209.857400000001

$ws.Range("H83").Value = 999.875
$ws.Range("I83").Value = 500
$ws.Range("J83").Value = 1071.2858
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 9641.572200000001
$ws.Range("M83").Value = 492
$ws.Range("N83").Value = -19625.5722

$ws.Range("H132").Value = 1742.2667
$ws.Range("I132").Value = 1748.976
$ws.Range("J132").Value = 1662.7142
$ws.Range("K132").Value = 5246.928
$ws.Range("L132").Value = 4988.142599999999
$ws.Range("M132").Value = -2716.928
$ws.Range("N132").Value = -10048.1426

$ws.Range("H137").Value = 2336.7917
$ws.Range("I137").Value = 2121.1428
$ws.Range("J137").Value = 2638.7
$ws.Range("K137").Value = 6363.428400000001
$ws.Range("L137").Value = 7916.099999999999
$ws.Range("M137").Value = -3813.428400000001
$ws.Range("N137").Value = -13016.1

$ws.Range("H138").Value = 3011.527
$ws.Range("I138").Value = 1437.2646
$ws.Range("K138").Value = 4311.793799999999
$ws.Range("M138").Value = 828.2062000000005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 19571.375
$ws.Range("I43").Value = 27447.5
$ws.Range("J43").Value = 16946
$ws.Range("K43").Value = 27447.5
$ws.Range("L43").Value = 16946
$ws.Range("M43").Value = -27134.5
$ws.Range("N43").Value = -17572

$ws.Range("H61").Value = 4014.0789
$ws.Range("I61").Value = 2715.8845
$ws.Range("J61").Value = 6826.8335
$ws.Range("K61").Value = 2715.8845
$ws.Range("L61").Value = 6826.8335
$ws.Range("M61").Value = -2503.8845
$ws.Range("N61").Value = -7250.8335

$ws.Range("H74").Value = 1939.2858
$ws.Range("I74").Value = 1616.9166
$ws.Range("K74").Value = 1616.9166
$ws.Range("M74").Value = -742.9166

$ws.Range("H77").Value = 1939.2858
$ws.Range("I77").Value = 1616.9166
$ws.Range("K77").Value = 8084.583000000001
$ws.Range("M77").Value = -3716.583000000001

$ws.Range("H132").Value = 3069.4358
$ws.Range("I132").Value = 2537.0312
$ws.Range("J132").Value = 5503.2856
$ws.Range("K132").Value = 7611.0936
$ws.Range("L132").Value = 16509.8568
$ws.Range("M132").Value = -5081.0936
$ws.Range("N132").Value = -21569.8568

$ws.Range("H136").Value = 4014.0789
$ws.Range("I136").Value = 2715.8845
$ws.Range("J136").Value = 6826.8335
$ws.Range("K136").Value = 8147.6535
$ws.Range("L136").Value = 20480.5005
$ws.Range("M136").Value = -5597.6535
$ws.Range("N136").Value = -25580.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2214.1458
$ws.Range("I134").Value = 2248.634
$ws.Range("J134").Value = 2012.1428
$ws.Range("K134").Value = 6745.902
$ws.Range("L134").Value = 6036.428400000001
$ws.Range("M134").Value = -4210.902
$ws.Range("N134").Value = -11106.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6030.6665
$ws.Range("I31").Value = 4958.76
$ws.Range("J31").Value = 8466.817999999999
$ws.Range("K31").Value = 4958.76
$ws.Range("L31").Value = 8466.817999999999
$ws.Range("M31").Value = -4663.76
$ws.Range("N31").Value = -9056.817999999999

$ws.Range("H34").Value = 6030.6665
$ws.Range("I34").Value = 4958.76
$ws.Range("J34").Value = 8466.817999999999
$ws.Range("K34").Value = 4958.76
$ws.Range("L34").Value = 8466.817999999999
$ws.Range("M34").Value = -4756.76
$ws.Range("N34").Value = -8870.817999999999

$ws.Range("H58").Value = 2680.1316
$ws.Range("I58").Value = 2592.52
$ws.Range("J58").Value = 2848.6155
$ws.Range("K58").Value = 2592.52
$ws.Range("L58").Value = 2848.6155
$ws.Range("M58").Value = -2389.52
$ws.Range("N58").Value = -3254.6155

$ws.Range("H132").Value = 1051.4534
$ws.Range("I132").Value = 952.69696
$ws.Range("J132").Value = 1775.6666
$ws.Range("K132").Value = 2858.09088
$ws.Range("L132").Value = 5326.9998
$ws.Range("M132").Value = -328.0908799999997
$ws.Range("N132").Value = -10386.9998

$ws.Range("H134").Value = 1300.7084
$ws.Range("I134").Value = 1239.8677
$ws.Range("K134").Value = 3719.6031
$ws.Range("M134").Value = -1184.6031

$ws.Range("H136").Value = 2680.1316
$ws.Range("I136").Value = 2592.52
$ws.Range("J136").Value = 2848.6155
$ws.Range("K136").Value = 7777.559999999999
$ws.Range("L136").Value = 8545.8465
$ws.Range("M136").Value = -5227.559999999999
$ws.Range("N136").Value = -13645.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2999
$ws.Range("J98").Value = 2999
$ws.Range("L98").Value = 8997
$ws.Range("N98").Value = -11993

$ws.Range("H113").Value = 2087.0625
$ws.Range("J113").Value = 2398.077
$ws.Range("L113").Value = 7194.231000000001
$ws.Range("N113").Value = -11534.231

$ws.Range("H137").Value = 3817.625
$ws.Range("I137").Value = 2497.0588
$ws.Range("J137").Value = 7024.7144
$ws.Range("K137").Value = 7491.176399999999
$ws.Range("L137").Value = 21074.1432
$ws.Range("M137").Value = -2391.176399999999
$ws.Range("N137").Value = -31274.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4093.1025
$ws.Range("I132").Value = 4124.919
$ws.Range("K132").Value = 12374.757
$ws.Range("M132").Value = -9844.757

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1543.6
$ws.Range("I7").Value = 1554.5
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 1554.5
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -1442.5
$ws.Range("N7").Value = -1724

$ws.Range("H59").Value = 9199
$ws.Range("J59").Value = 9199
$ws.Range("L59").Value = 9199
$ws.Range("N59").Value = -10507

$ws.Range("H100").Value = 3999.6667
$ws.Range("I100").Value = 3299.6
$ws.Range("J100").Value = 7500
$ws.Range("K100").Value = 3299.6
$ws.Range("L100").Value = 7500
$ws.Range("M100").Value = -2758.6
$ws.Range("N100").Value = -8582

$ws.Range("H126").Value = 1543.6
$ws.Range("I126").Value = 1554.5
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 4663.5
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -2193.5
$ws.Range("N126").Value = -9440

$ws.Range("H132").Value = 18399.666
$ws.Range("I132").Value = 19033
$ws.Range("J132").Value = 17133
$ws.Range("K132").Value = 57099
$ws.Range("L132").Value = 51399
$ws.Range("M132").Value = -54569
$ws.Range("N132").Value = -56459

$ws.Range("H136").Value = 1880.1951
$ws.Range("I136").Value = 1910.5135
$ws.Range("J136").Value = 1599.75
$ws.Range("K136").Value = 5731.5405
$ws.Range("L136").Value = 4799.25
$ws.Range("M136").Value = -3181.5405
$ws.Range("N136").Value = -9899.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 25250
$ws.Range("J114").Value = 25250
$ws.Range("L114").Value = 25250
$ws.Range("N114").Value = -33928

$ws.Range("H132").Value = 2623.4814
$ws.Range("I132").Value = 2508.7083
$ws.Range("J132").Value = 3541.6667
$ws.Range("K132").Value = 7526.124899999999
$ws.Range("L132").Value = 10625.0001
$ws.Range("M132").Value = -4996.124899999999
$ws.Range("N132").Value = -15685.0001
